# Update the "data_monitoreo_la_peñita" sheet: the empadronador totals table
# was re-sorted (descending by total_registros) and some counts were bumped
# up for newly-added registrations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering / values for rows 2-11 (rows 12-14 remain unchanged).
$data = @(
    @("GARAVITO LEON IVONNE LISSETH", 115),
    @("TIMOTEO BAYONA SHARYN LISSETH", 114),
    @("ZAPATA ZETA ROSA ARACELI", 112),
    @("PEREZ VEGA ANA YSABEL", 110),
    @("NIÑO GUERRERO ANYELA MELINA", 91),
    @("PANTA MONZON SHIRLEY MARIBEL", 89),
    @("VALLE SILVA SUTMMER ORFELINDA", 82),
    @("TIZON NUÑEZ FRESIA YAMILI", 80),
    @("CASTRO JUAREZ MARIA ISABEL", 79),
    @("MORENO PALACIOS DAMARIS VANESA", 74)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
